$d = $word.ActiveDocument

# Locate the paragraph that currently contains "To," followed by a tab and
# the bold "Notice u/s 94 BNSS, 2023" run.
$toParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -eq "To,`tNotice u/s 94 BNSS, 2023`r") {
        $toParaIndex = $i
        break
    }
}

if ($toParaIndex -eq -1) {
    Write-Host "Could not locate the target paragraph"
} else {
    $toPara = $d.Paragraphs.Item($toParaIndex)

    # 1) Split a brand-new empty paragraph in right before the "To," paragraph.
    $toPara.Range.InsertParagraphBefore()

    # The freshly inserted (still empty) paragraph now occupies the old index.
    $noticePara = $d.Paragraphs.Item($toParaIndex)
    $noticeRange = $noticePara.Range

    # 2) Replace that empty paragraph's contents with the centered, bold +
    #    underlined "Notice u/s 94 BNSS, 2023" heading, with no other
    #    paragraph formatting carried over from the old paragraph.
    $noticeXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' +
      '<w:p>' +
      '<w:pPr><w:jc w:val="center"/></w:pPr>' +
      '<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Notice u/s 94 BNSS, 2023</w:t></w:r>' +
      '</w:p>' +
      '</w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'
    $noticeRange.InsertXML($noticeXml, "Replace")

    # 3) Clean up the original "To," paragraph (now shifted one slot down):
    #    drop the centered tab stop, the tab character, and the bold
    #    "Notice u/s 94 BNSS, 2023" run, leaving a plain "To," line that is
    #    explicitly left aligned.
    $toPara = $d.Paragraphs.Item($toParaIndex + 1)
    $toPara.Range.ParagraphFormat.TabStops.ClearAll()

    $toParaStart = $toPara.Range.Start
    $toParaEnd = $toPara.Range.End
    $tailRange = $d.Range($toParaStart + 3, $toParaEnd - 1)
    $tailRange.Delete()

    $toPara = $d.Paragraphs.Item($toParaIndex + 1)
    $toPara.Alignment = 0

    Write-Host "Notice heading moved above 'To,' and reformatted."
}
